$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.438.33"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.725.13"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4916"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2616"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "1.720.51"
$ws.Range("E10").Value = "  -0.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06994"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.565"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5996"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "26.441.51"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007165"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.82%  "
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").Value = "1.941.76"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.480"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.598"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.163"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.399"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "107.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.957"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07962"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.674"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04521"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6264"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9279"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.957"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.34%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.389"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9996"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01488"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.346"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.01%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3848"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.762"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.48%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1166"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05367"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.699"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.229"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.85%  "
